# Operator-interface BOM: consolidate the Red/Green LED line items into a
# single Cool White LED line (with an Amber alternative note), per commit
# "Updating LEDs in BOM to white / amber".
#
# Before:
#   Row 23: Qty 8, (blank Value), LED5MM, LED1..LED6, LED7/LED9,  C503B-GAN...-ND (Green)
#   Row 24: Qty 2, (blank Value), LED5MM, (blank refs), LED8/LED10, C503B-RCN...-ND (Red)
# After:
#   Row 23: Qty 10, Cool White, LED5MM, LED1..LED6, LED7/LED8/LED9/LED10,
#           C503C-WAN-CBBDB231-ND, note "Alternative: C503C-ACN-CYCZA342CT-ND (Amber)"
#   (old row 24 is gone; every row below shifts up by one)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Red" LED row entirely -- its content is being folded into
# row 23, and deleting it shifts rows 25-29 up to 24-28 (matching the new
# dimension A1:H28).
$ws.Rows.Item(24).Delete()

# Update row 23 in place to describe the combined white LED part.
$ws.Cells.Item(23, 1).Value = 10
$ws.Cells.Item(23, 5).Value = "LED7, LED8, LED9, LED10"
$ws.Cells.Item(23, 8).Value = "Alternative: C503C-ACN-CYCZA342CT-ND (Amber)"
$ws.Cells.Item(23, 2).Value = "Cool White"
$ws.Cells.Item(23, 6).Value = "C503C-WAN-CBBDB231-ND"

# Leave the selection where the author's last edit was.
$ws.Range("F23").Select()
